$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Title placeholder ("ctrTitle") - give it an explicit position/size
# ---------------------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Left   = 54.0
$title.Top    = 46.24953
$title.Width  = 612.0
$title.Height = 115.75

# ---------------------------------------------------------------------
# 2) Subtitle placeholder ("subTitle") - explicit position/size, and
#    merge the "Darci" + " " runs into a single "Darci " run (matches
#    what PowerPoint does when you retype across run boundaries).
# ---------------------------------------------------------------------
$subtitle = $s.Shapes.Item(2)
$subtitle.Left   = 108.0
$subtitle.Top    = 201.3745
$subtitle.Width  = 504.0
$subtitle.Height = 138.0

$subTr = $subtitle.TextFrame.TextRange
# Remove the original "Darci" (chars 1-5), then type it again right
# before the space so it inherits the space run's (non-misspelled)
# formatting and merges into one run "Darci ".
$subTr.Characters(1, 5).Text = ""
$subTr.Characters(1, 1).InsertBefore("Darci")

# ---------------------------------------------------------------------
# 3) New shape: copyright / license notice rectangle
# ---------------------------------------------------------------------
$rect = $s.Shapes.AddShape(1, 27.5625297, 385.8375, 664.875, 116.32504)
$rect.Name = "Rectangle 4"
$rect.TextFrame.WordWrap = -1
$rect.TextFrame.AutoSize = 1

$rTr = $rect.TextFrame.TextRange
$full = "Copyright 2018 Darci Burdge and Stoney Jackson SOME RIGHTS RESERVED" + [char]13 + [char]13 + "This work is licensed under the Creative Commons Attribution-ShareAlike 4.0 International License. To view a copy of this license, visit http://creativecommons.org/licenses/by-sa/4.0/ ."
$rTr.Text = $full

# Re-assign (no-op content change) the spans that need to land in their
# own run so formatting / run-boundaries match the source edit.
function Set-RunSpan($textRange, $needle, $startFrom) {
    $idx = $textRange.Text.IndexOf($needle, $startFrom - 1)
    $pos = $idx + 1
    $len = $needle.Length
    $textRange.Characters($pos, $len).Text = $needle
    return $pos + $len
}

$cursor = 1
$cursor = Set-RunSpan $rTr "Copyright 2018 Darci " $cursor
$cursor = Set-RunSpan $rTr "Burdge" $cursor
$cursor = Set-RunSpan $rTr " and Stoney Jackson SOME RIGHTS RESERVED" $cursor

$p2 = $rTr.Paragraphs(3, 1)
$base = $p2.Start
$cursor = $base
$cursor = Set-RunSpan $rTr "This work is licensed under the Creative Commons Attribution-" $cursor
$cursor = Set-RunSpan $rTr "ShareAlike" $cursor
$cursor = Set-RunSpan $rTr " 4.0 International License. To view a copy of this license, visit http://" $cursor
$cursor = Set-RunSpan $rTr "creativecommons.org" $cursor
$cursor = Set-RunSpan $rTr "/licenses/by-" $cursor
$cursor = Set-RunSpan $rTr "sa" $cursor
$cursor = Set-RunSpan $rTr "/4.0/ ." $cursor

Write-Output "edit complete"
